$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 9).Value = 'b'
$ws.Cells.Item(2, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(3, 9).Value = 'sd'
$ws.Cells.Item(3, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(11, 9).Value = 'sv'
$ws.Cells.Item(11, 10).Value = 'Statement-opinion'
$ws.Cells.Item(15, 9).Value = 'sd'
$ws.Cells.Item(15, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(34, 9).Value = 'sd'
$ws.Cells.Item(34, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(41, 9).Value = 'b'
$ws.Cells.Item(41, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(54, 9).Value = 'sd'
$ws.Cells.Item(54, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(56, 9).Value = 'ba'
$ws.Cells.Item(56, 10).Value = 'Appreciation'
$ws.Cells.Item(63, 9).Value = 'aa'
$ws.Cells.Item(63, 10).Value = 'Agree/Accept'
$ws.Cells.Item(70, 9).Value = 'aa'
$ws.Cells.Item(70, 10).Value = 'Agree/Accept'
$ws.Cells.Item(71, 9).Value = 'sd'
$ws.Cells.Item(71, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(75, 9).Value = 'ba'
$ws.Cells.Item(75, 10).Value = 'Appreciation'
$ws.Cells.Item(77, 9).Value = 'ba'
$ws.Cells.Item(77, 10).Value = 'Appreciation'
$ws.Cells.Item(79, 9).Value = 'sv'
$ws.Cells.Item(79, 10).Value = 'Statement-opinion'
$ws.Cells.Item(85, 9).Value = 'sd'
$ws.Cells.Item(85, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(101, 9).Value = 'sd'
$ws.Cells.Item(101, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(128, 9).Value = 'sv'
$ws.Cells.Item(128, 10).Value = 'Statement-opinion'
$ws.Cells.Item(134, 9).Value = 'sv'
$ws.Cells.Item(134, 10).Value = 'Statement-opinion'
$ws.Cells.Item(138, 9).Value = 'sd'
$ws.Cells.Item(138, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(145, 9).Value = 'b'
$ws.Cells.Item(145, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(146, 9).Value = 'aa'
$ws.Cells.Item(146, 10).Value = 'Agree/Accept'
$ws.Cells.Item(149, 9).Value = 'sv'
$ws.Cells.Item(149, 10).Value = 'Statement-opinion'
$ws.Cells.Item(155, 9).Value = 'aa'
$ws.Cells.Item(155, 10).Value = 'Agree/Accept'
$ws.Cells.Item(164, 9).Value = 'ba'
$ws.Cells.Item(164, 10).Value = 'Appreciation'
$ws.Cells.Item(177, 9).Value = 'ba'
$ws.Cells.Item(177, 10).Value = 'Appreciation'
$ws.Cells.Item(180, 9).Value = 'aa'
$ws.Cells.Item(180, 10).Value = 'Agree/Accept'
$ws.Cells.Item(185, 9).Value = 'sd'
$ws.Cells.Item(185, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(191, 9).Value = 'aa'
$ws.Cells.Item(191, 10).Value = 'Agree/Accept'
$ws.Cells.Item(195, 9).Value = '%'
$ws.Cells.Item(195, 10).Value = 'Uninterpretable'
$ws.Cells.Item(207, 9).Value = 'sd'
$ws.Cells.Item(207, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(208, 9).Value = 'sd'
$ws.Cells.Item(208, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(209, 9).Value = 'sd'
$ws.Cells.Item(209, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(214, 9).Value = 'aa'
$ws.Cells.Item(214, 10).Value = 'Agree/Accept'
$ws.Cells.Item(217, 9).Value = 'sd'
$ws.Cells.Item(217, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(221, 9).Value = '%'
$ws.Cells.Item(221, 10).Value = 'Uninterpretable'
$ws.Cells.Item(241, 9).Value = '%'
$ws.Cells.Item(241, 10).Value = 'Uninterpretable'
$ws.Cells.Item(242, 9).Value = 'b'
$ws.Cells.Item(242, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(243, 9).Value = 'sd'
$ws.Cells.Item(243, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(245, 9).Value = 'sd'
$ws.Cells.Item(245, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(253, 9).Value = '%'
$ws.Cells.Item(253, 10).Value = 'Uninterpretable'
$ws.Cells.Item(254, 9).Value = '%'
$ws.Cells.Item(254, 10).Value = 'Uninterpretable'
$ws.Cells.Item(260, 9).Value = 'aa'
$ws.Cells.Item(260, 10).Value = 'Agree/Accept'
$ws.Cells.Item(266, 9).Value = '%'
$ws.Cells.Item(266, 10).Value = 'Uninterpretable'
$ws.Cells.Item(275, 9).Value = 'sd'
$ws.Cells.Item(275, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(276, 9).Value = 'aa'
$ws.Cells.Item(276, 10).Value = 'Agree/Accept'
$ws.Cells.Item(278, 9).Value = 'aa'
$ws.Cells.Item(278, 10).Value = 'Agree/Accept'
$ws.Cells.Item(284, 9).Value = 'sd'
$ws.Cells.Item(284, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(298, 9).Value = 'sd'
$ws.Cells.Item(298, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(299, 9).Value = 'sd'
$ws.Cells.Item(299, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(300, 9).Value = 'sd'
$ws.Cells.Item(300, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(302, 9).Value = 'sd'
$ws.Cells.Item(302, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(303, 9).Value = 'sd'
$ws.Cells.Item(303, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(313, 9).Value = 'sd'
$ws.Cells.Item(313, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(314, 9).Value = 'sd'
$ws.Cells.Item(314, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(316, 9).Value = 'ba'
$ws.Cells.Item(316, 10).Value = 'Appreciation'
$ws.Cells.Item(317, 9).Value = 'sd'
$ws.Cells.Item(317, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(319, 9).Value = 'qy'
$ws.Cells.Item(319, 10).Value = 'Yes-No-Question'
$ws.Cells.Item(320, 9).Value = 'sd'
$ws.Cells.Item(320, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(329, 9).Value = 'sd'
$ws.Cells.Item(329, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(331, 9).Value = 'sv'
$ws.Cells.Item(331, 10).Value = 'Statement-opinion'
$ws.Cells.Item(333, 9).Value = 'ba'
$ws.Cells.Item(333, 10).Value = 'Appreciation'
$ws.Cells.Item(338, 9).Value = 'ba'
$ws.Cells.Item(338, 10).Value = 'Appreciation'
$ws.Cells.Item(339, 9).Value = 'b'
$ws.Cells.Item(339, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(346, 9).Value = 'sd'
$ws.Cells.Item(346, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(357, 9).Value = 'aa'
$ws.Cells.Item(357, 10).Value = 'Agree/Accept'
$ws.Cells.Item(362, 9).Value = '%'
$ws.Cells.Item(362, 10).Value = 'Uninterpretable'
$ws.Cells.Item(363, 9).Value = '%'
$ws.Cells.Item(363, 10).Value = 'Uninterpretable'
$ws.Cells.Item(369, 9).Value = 'b'
$ws.Cells.Item(369, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(370, 9).Value = 'ba'
$ws.Cells.Item(370, 10).Value = 'Appreciation'
$ws.Cells.Item(371, 9).Value = 'aa'
$ws.Cells.Item(371, 10).Value = 'Agree/Accept'
$ws.Cells.Item(374, 9).Value = 'sd'
$ws.Cells.Item(374, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(375, 9).Value = 'sv'
$ws.Cells.Item(375, 10).Value = 'Statement-opinion'
$ws.Cells.Item(377, 9).Value = 'aa'
$ws.Cells.Item(377, 10).Value = 'Agree/Accept'
$ws.Cells.Item(391, 9).Value = 'aa'
$ws.Cells.Item(391, 10).Value = 'Agree/Accept'
$ws.Cells.Item(414, 9).Value = 'sd'
$ws.Cells.Item(414, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(417, 9).Value = '%'
$ws.Cells.Item(417, 10).Value = 'Uninterpretable'
$ws.Cells.Item(418, 9).Value = '%'
$ws.Cells.Item(418, 10).Value = 'Uninterpretable'
$ws.Cells.Item(425, 9).Value = 'aa'
$ws.Cells.Item(425, 10).Value = 'Agree/Accept'
$ws.Cells.Item(430, 9).Value = 'ba'
$ws.Cells.Item(430, 10).Value = 'Appreciation'
$ws.Cells.Item(432, 9).Value = 'b'
$ws.Cells.Item(432, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(439, 9).Value = 'sd'
$ws.Cells.Item(439, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(441, 9).Value = 'aa'
$ws.Cells.Item(441, 10).Value = 'Agree/Accept'
$ws.Cells.Item(456, 9).Value = 'sv'
$ws.Cells.Item(456, 10).Value = 'Statement-opinion'
$ws.Cells.Item(460, 9).Value = 'sv'
$ws.Cells.Item(460, 10).Value = 'Statement-opinion'
$ws.Cells.Item(465, 9).Value = 'ba'
$ws.Cells.Item(465, 10).Value = 'Appreciation'
$ws.Cells.Item(475, 9).Value = 'ba'
$ws.Cells.Item(475, 10).Value = 'Appreciation'
$ws.Cells.Item(476, 9).Value = 'sd'
$ws.Cells.Item(476, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(479, 9).Value = 'b'
$ws.Cells.Item(479, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(480, 9).Value = 'sd'
$ws.Cells.Item(480, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(482, 9).Value = 'b'
$ws.Cells.Item(482, 10).Value = 'Acknowledge (Backchannel)'
